$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Unit Number Info")

# Clear the header labels in A1/B1 (keep formatting/style), and remove the
# data row (row 2) entirely so coordinates shift as expected.
$ws.Range("A1:B1").ClearContents()
$ws.Range("A2:B2").ClearContents()
